# Fixing Concurrent Modification Exception
# Enable (set isEnabled = TRUE) for a set of ingredients on the AllIngredients sheet,
# and update the frozen-pane/selection view state to reflect scrolled position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllIngredients")

# Rows (by ingredient name) whose "isEnabled" (column B) flag should become TRUE
$namesToEnable = @(
    "Basileus",
    "Blood Ore",
    "Common Vitus Juice",
    "Compact Horn",
    "Electrum",
    "Gamun Oil",
    "Gold",
    "Green Jamburra",
    "Green Jamburra Juice",
    "Kyanite",
    "Pirum",
    "Pirum Juice",
    "Waterstone",
    "White Bear Carcass"
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($namesToEnable -contains $name) {
        $ws.Cells.Item($r, 2).Value = $true
    }
}

# Update the view state: scroll so row 11 is the top-left cell of the frozen pane,
# and move the active selection to C16.
$ws.Activate()
$ws.Range("A11").Select()
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("C16").Select()
